# bercot/time_sheet.xlsx - "error messages show up in HTML, some basic
# prettying up completed"
#
# The sheet tracked one time-in/time-out pair per day (columns C/D, stored
# as fractional-day clock times). This switches rows 6-7 over to a single
# "hours worked" figure typed straight into column B (like the rows above
# them already do), adds a new day (row 8) that still uses the old
# time-in/time-out style in column C, and appends a running total of the
# hours column at row 18. A bit of column widening and cursor repositioning
# rounds out the "prettying up".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: drop the old C6/D6 clock-time pair, enter hours directly in B6
$ws.Range("B6").Value = 1.33
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()

# --- Row 7: drop the old C7 clock-time value, enter hours directly in B7
$ws.Range("B7").Value = 1.2
$ws.Range("C7").ClearContents()

# --- Row 8 (new): another day, still logged the old way (date + C time)
# Copy formatting from the row above first so the new cells pick up the
# same number formats (date / time-of-day) without minting new styles.
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(8, 1).Value = 42930

$ws.Cells.Item(7, 3).Copy()
$ws.Cells.Item(8, 3).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(8, 3).Value = 0.149305555555556

$excel.CutCopyMode = $false

# --- Row 18 (new): running total of the hours column
$ws.Range("B18").Formula = "=SUM(B2:B15)"

# --- Widen the columns a touch
$ws.Range("A:B").ColumnWidth = 7.6
$ws.Columns.Item(3).ColumnWidth = 9.92
$ws.Columns.Item(4).ColumnWidth = 10.1

# --- Leave the cursor where the author left it
$ws.Range("B9").Select()
